$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Q0) ---
$ws.Cells.Item(2,2).Value = -0.2086161574976735
$ws.Cells.Item(2,3).Value = 0.4831534245497861
$ws.Cells.Item(2,4).Value = 0.3849075325051517
$ws.Cells.Item(2,5).Value = 0.6204091653942192
$ws.Cells.Item(2,6).Value = 0.6063393083278098

# --- Row 3 (Q1) ---
$ws.Cells.Item(3,2).Value = -0.130351666236421
$ws.Cells.Item(3,3).Value = 0.4021978896510516
$ws.Cells.Item(3,4).Value = 0.2468797842177299
$ws.Cells.Item(3,5).Value = 0.4968699872378386
$ws.Cells.Item(3,6).Value = 0.4990446669431289

# --- Row 4 (Q2) ---
$ws.Cells.Item(4,2).Value = -0.172168584493037
$ws.Cells.Item(4,3).Value = 0.4469611009306131
$ws.Cells.Item(4,4).Value = 0.3068024783553724
$ws.Cells.Item(4,5).Value = 0.5538975341661782
$ws.Cells.Item(4,6).Value = 0.5498698591839245

# --- Row 5 (Q3) ---
$ws.Cells.Item(5,2).Value = -0.2012538732070167
$ws.Cells.Item(5,3).Value = 0.3706701098662007
$ws.Cells.Item(5,4).Value = 0.2450590527804206
$ws.Cells.Item(5,5).Value = 0.4950343955528955
$ws.Cells.Item(5,6).Value = 0.4743537966850842

# --- Row 6 (Q4) ---
$ws.Cells.Item(6,2).Value = -0.1962591978527491
$ws.Cells.Item(6,3).Value = 0.5005154228308657
$ws.Cells.Item(6,4).Value = 0.3060069273344895
$ws.Cells.Item(6,5).Value = 0.5531789288598126
$ws.Cells.Item(6,6).Value = 0.5451699577936784

# --- Row 7 (Q5) ---
$ws.Cells.Item(7,2).Value = -0.2173851562319015
$ws.Cells.Item(7,3).Value = 0.5193912335752287
$ws.Cells.Item(7,4).Value = 0.3417199285034996
$ws.Cells.Item(7,5).Value = 0.5845681555674236
$ws.Cells.Item(7,6).Value = 0.575561964646486
$ws.Cells.Item(7,7).Value = 9

# --- Row 8 (Q6) ---
$ws.Cells.Item(8,2).Value = -0.1465976859772833
$ws.Cells.Item(8,3).Value = 0.6786390806747216
$ws.Cells.Item(8,4).Value = 0.5659560571029979
$ws.Cells.Item(8,5).Value = 0.7523005098383211
$ws.Cells.Item(8,6).Value = 0.8083057655880753
$ws.Cells.Item(8,7).Value = 6

# --- Row 9 (Q7) ---
$ws.Cells.Item(9,2).Value = -0.7439835285872217
$ws.Cells.Item(9,3).Value = 0.7439835285872217
$ws.Cells.Item(9,4).Value = 0.7169109747095641
$ws.Cells.Item(9,5).Value = 0.8467059552817401
$ws.Cells.Item(9,6).Value = 0.4950749699295109
$ws.Cells.Item(9,7).Value = 3

# --- Row 10 (NEW, Q8) ---
$ws.Cells.Item(9,1).Copy()
$ws.Cells.Item(10,1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Cells.Item(10,1).Value = "Q8"
$ws.Cells.Item(10,2).Value = -0.8456720827280808
$ws.Cells.Item(10,3).Value = 0.8456720827280808
$ws.Cells.Item(10,4).Value = 0.7151612715056499
$ws.Cells.Item(10,5).Value = 0.8456720827280808
$ws.Cells.Item(10,7).Value = 1
